$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update capital structure database values for rows 2 and 3 (both rows share identical values)
foreach ($r in 2,3) {
    $ws.Range("G$r").Value  = -0.1849489795918367
    $ws.Range("H$r").Value  = -0.1849489795918367
    $ws.Range("I$r").Value  = -0.1658163265306123
    $ws.Range("J$r").Value  = -0.1658163265306123
    $ws.Range("K$r").Value  = -1.96
    $ws.Range("L$r").Value  = -0.25
    $ws.Range("M$r").Value  = 0.12
    $ws.Range("N$r").Value  = 0.003658536585365854
    $ws.Range("O$r").Value  = -0.06122448979591837
    $ws.Range("S$r").Value  = 0.12
    $ws.Range("T$r").Value  = 1
    $ws.Range("U$r").Value  = 0.163
    $ws.Range("V$r").Value  = 0.004969512195121952
    $ws.Range("W$r").Value  = -1.248407643312102
    $ws.Range("X$r").Value  = 0.06026819411833331
    $ws.Range("Y$r").Value  = -1.308675837430435
    $ws.Range("Z$r").Value  = 9.987261146496815
    $ws.Range("AA$r").Value = -1.656050955414013
    $ws.Range("AB$r").Value = 0.06026819411833331
    $ws.Range("AC$r").Value = -1.716319149532346
    $ws.Range("AG$r").Value = -0.163
    $ws.Range("AJ$r").Value = -0.004994331586849281
    $ws.Range("AK$r").Value = -0.2620578778135048
    $ws.Range("AL$r").Value = 0.02
    $ws.Range("AM$r").Value = 0.02
    $ws.Range("AN$r").Value = -0
    $ws.Range("AO$r").Value = -65
    $ws.Range("AP$r").Value = 0.12734375
    $ws.Range("AQ$r").Value = -65
}
